# "Fruta / hortaliza, semanal" — insert a new weekly price observation row
# for Coliflor @ Vega Modelo de Temuco. This pushes the existing rows
# 175..228 down to 176..229 (Excel's normal Insert-shift-down behaviour)
# and the brand-new data point lands in row 175.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 175:228 down by one row, like picking Insert on the row header.
$ws.Rows("175:175").Insert()

# Populate the newly inserted row 175 with the new observation.
$ws.Cells.Item(175, 1).Value = 10
$ws.Cells.Item(175, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(175, 3).Value = "La Araucanía"
$ws.Cells.Item(175, 4).Value = 44463
$ws.Cells.Item(175, 5).Value = 9
$ws.Cells.Item(175, 6).Value = 100112008
$ws.Cells.Item(175, 7).Value = "Coliflor"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 300
$ws.Cells.Item(175, 11).Value = 800
$ws.Cells.Item(175, 12).Value = 800
$ws.Cells.Item(175, 13).Value = 800
$ws.Cells.Item(175, 14).Value = "$/unidad"
$ws.Cells.Item(175, 15).Value = "Región Metropolitana"
$ws.Cells.Item(175, 16).Value = 800
$ws.Cells.Item(175, 17).Value = 1
$ws.Cells.Item(175, 18).Value = "Hortaliza"
